# Apply edits described in the commit:
# - Adjust the "N", "O" and "S" values for instance02 (row 5) to account for the
#   time elapsed since the last improvement when renewing the population.
# - Widen column A so the longer instance labels fit ("tabulist" note column).
# - Move the active selection to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Update the raw benchmark inputs for instance02 (row 5)
$ws.Range("N5").Value = 42.89
$ws.Range("O5").Value = 40.79
$ws.Range("S5").Value = 11.32

# Widen column A to fit the updated content (target stored width = 10.5)
$ws.Columns.Item(1).ColumnWidth = 9.666666666666666

# Recalculate dependent formulas (AVERAGE/Gap% cells) explicitly
$excel.Calculate()

# Move/resize the active cell selection
$ws.Range("A8").Select()
